$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "子供オペレーターたちの授業を終えたギターノは、次の授業を行うためにやってきたクロワッサンと会話をする。クロワッサンは商売で問題に直面しているらしい。ギターノは報酬を受け取って、彼女のために占いをするのだった。`n"
$ws.Range("C2").Value = "Gitano finishes her lesson with the young Operators when she meets Croissant, who's coming to teach the next. Croissant seems to have run into some issues with her business, and after receiving compensation, Gitano divines Croissant's fortune for her.`n"
